# Delete row 484 ("戦士" entry) from the sheet.
# This removes the whole row and shifts all subsequent rows up by one,
# which matches the diff (rows 485-573 become 484-572, dimension C573 -> C572).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(484).Delete()
